# Refresh crypto price/volume snapshot (mirrors the upstream data-refresh commit).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '51.047.91'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.32%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.947.11'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.17%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '378.78'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.66%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '104.47'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.01%  '
$ws.Range('E7').Value = '  +1.15%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.593'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.73%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.97'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.81%  '
$ws.Range('E11').Value = '  +0.68%  '
$ws.Range('E12').Value = '  +1.27%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.36'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.80%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.408.16'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.87%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.48'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.87%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.932.99'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.84%  '
$ws.Range('E17').Value = '  +3.78%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '51.041.18'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.43%  '
$ws.Range('E19').Value = '  +2.47%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.35'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.87%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.87'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.24%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0₃0958'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.53%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '69.39'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.44%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '261.00'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.61%  '
$ws.Range('E25').Value = '  +5.47%  '
$ws.Range('B26').Value = 'Filecoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.61'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +10.43%  '
$ws.Range('B27').Value = 'RenderToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.28'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +22.85%  '
$ws.Range('E28').Value = '  +1.19%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.112'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +9.43%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '25.80'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.27%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '9.82'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.33%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '34.51'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.92%  '
$ws.Range('E34').Value = '  -2.14%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '50.87'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.24%  '
$ws.Range('E36').Value = '  +8.33%  '
$ws.Range('E37').Value = '  +0.07%  '
$ws.Range('E38').Value = '  +1.65%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '17.20'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.41%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.57'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.07%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.84'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.46%  '
$ws.Range('E42').Value = '  +3.25%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '122.26'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.03%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '21.95'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.35%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.288'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +23.76%  '
$ws.Range('E46').Value = '  -0.93%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.39'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.37%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.032.03'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.66%  '
$ws.Range('E49').Value = '  +3.02%  '
$ws.Range('E50').Value = '  +11.94%  '
$ws.Range('E51').Value = '  +2.75%  '
